$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 10 de Septiembre de 2020 a las 22:33"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 6576582
$ws.Range("C4").Value = 27107
$ws.Range("D4").Value = 3867321
$ws.Range("E4").Value = 2513251
$ws.Range("G4").Value = 772
$ws.Range("H4").Value = 196010

# Row 5: India
$ws.Range("D5").Value = 3539983
$ws.Range("E5").Value = 943438

# Row 24: Alemania
$ws.Range("B24").Value = 258107
$ws.Range("C24").Value = 1758
$ws.Range("E24").Value = 16788
$ws.Range("G24").Value = 9
$ws.Range("H24").Value = 9419

# Row 28: Israel
$ws.Range("B28").Value = 145526
$ws.Range("C28").Value = 4429
$ws.Range("D28").Value = 109942
$ws.Range("E28").Value = 34507
$ws.Range("G28").Value = 23
$ws.Range("H28").Value = 1077

# Row 43 & 44: Marruecos / Guatemala swap (names + new values)
$ws.Range("A43").Value = "Guatemala"
$ws.Range("B43").Value = 80306
$ws.Range("C43").Value = 684
$ws.Range("D43").Value = 68927
$ws.Range("E43").Value = 8461
$ws.Range("G43").Value = 21
$ws.Range("H43").Value = 2918

$ws.Range("A44").Value = "Marruecos"
$ws.Range("B44").Value = 79767
$ws.Range("C44").Value = 1889
$ws.Range("D44").Value = 61850
$ws.Range("E44").Value = 16426
$ws.Range("G44").Value = 38
$ws.Range("H44").Value = 1491

# Row 57: Costa Rica
$ws.Range("B57").Value = 52549
$ws.Range("C57").Value = 1325
$ws.Range("D57").Value = 20322
$ws.Range("E57").Value = 31660
$ws.Range("G57").Value = 24
$ws.Range("H57").Value = 567

# Row 98: Guayana Francesa
$ws.Range("D98").Value = 9040
$ws.Range("E98").Value = 359

# Row 104: Haiti
$ws.Range("B104").Value = 8429
$ws.Range("C104").Value = 45
$ws.Range("E104").Value = 2223
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 215

# Row 106: Mauritania
$ws.Range("B106").Value = 7222
$ws.Range("C106").Value = 31
$ws.Range("D106").Value = 6758
$ws.Range("E106").Value = 303

# Row 109 & 110: Malaui / Tunez swap
$ws.Range("A109").Value = "Tunez"
$ws.Range("B109").Value = 5882
$ws.Range("C109").Value = 465
$ws.Range("D109").Value = 1862
$ws.Range("E109").Value = 3921
$ws.Range("G109").Value = 3
$ws.Range("H109").Value = 99

$ws.Range("A110").Value = "Malaui"
$ws.Range("B110").Value = 5653
$ws.Range("C110").Value = 0
$ws.Range("D110").Value = 3630
$ws.Range("E110").Value = 1847
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 176

# Row 111: Republica de Yibuti
$ws.Range("B111").Value = 5394
$ws.Range("C111").Value = 3
$ws.Range("E111").Value = 6

# Row 150: Sierra Leona
$ws.Range("B150").Value = 2069
$ws.Range("C150").Value = 2
$ws.Range("E150").Value = 375

# Row 151: Yemen
$ws.Range("B151").Value = 2003
$ws.Range("C151").Value = 4
$ws.Range("D151").Value = 1211
$ws.Range("E151").Value = 212
$ws.Range("G151").Value = 4
$ws.Range("H151").Value = 580

# Row 156: Togo
$ws.Range("B156").Value = 1537
$ws.Range("C156").Value = 9
$ws.Range("D156").Value = 1145
$ws.Range("E156").Value = 355
$ws.Range("G156").Value = 1
$ws.Range("H156").Value = 37

# Row 166: Republica del Chad
$ws.Range("B166").Value = 1051
$ws.Range("C166").Value = 3
$ws.Range("D166").Value = 932
$ws.Range("E166").Value = 40

# Row 214 & 215: Montserrat / Islas Malvinas swap
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1
